$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 160; $r -le 186; $r++) {
    $ws.Cells.Item($r, 2).Value = 45993
}

$ws.Range("H180").Select() | Out-Null
